$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date (2021-09-26 -> 2021-09-27)
$ws.Name = "Through 2021-09-27"

# Update the "September (through ...)" row label to match the new date
$ws.Range("A10").Value = "September (through 09-27)"

# Update September row (row 10) counts for the new data pull
$ws.Range("B10").Value = 28
$ws.Range("C10").Value = 41
$ws.Range("D10").Value = 69
$ws.Range("F10").Value = 65
$ws.Range("G10").Value = 103
$ws.Range("H10").Value = 165

# Update Total row (row 11) counts to reflect the new September totals
$ws.Range("B11").Value = 222
$ws.Range("C11").Value = 422
$ws.Range("D11").Value = 620
$ws.Range("F11").Value = 414
$ws.Range("G11").Value = 887
$ws.Range("H11").Value = 1235
